# Adds the "Admin block/unblock wallet" user story estimation rows to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 57 - new User story (column A only)
$ws.Cells.Item(57, 1).Value = "Admin moze da blokira/odblokira odredjeni novcanik/nalog"

# Row 58 - task + estimated time
$ws.Cells.Item(58, 2).Value = "Dodavanje admin passworda u appsetttigns"
$ws.Cells.Item(58, 3).Value = 1

# Row 59 - task + estimated time
$ws.Cells.Item(59, 2).Value = "Dodavanje BlockWallet i UnblockWallet metoda u WalletService"
$ws.Cells.Item(59, 3).Value = 15

# Row 60 - task + estimated time
$ws.Cells.Item(60, 2).Value = "Implementacija testova za BlockWallet I UnblockWallet"
$ws.Cells.Item(60, 3).Value = 20

# Row 61 - task + estimated time
$ws.Cells.Item(61, 2).Value = "Dodavanje rute za blokiranje na WalletController"
$ws.Cells.Item(61, 3).Value = 10

# Row 62 - task + estimated time
$ws.Cells.Item(62, 2).Value = "Dodavanje stranice za blokiranje u MVC aplikaciju"
$ws.Cells.Item(62, 3).Value = 20

# Match the saved selection/active cell from the diff
$ws.Range("D63").Select()
